# Update crypto price/volume figures per the Tue Aug  1 09:00:56 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.955.43"
$ws.Range("E2").Value = "  -1.55%  "
$ws.Range("D3").Value = "1.834.34"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.16"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6900"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07702"
$ws.Range("E8").Value = "  -3.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3056"
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.58"
$ws.Range("E10").Value = "  -4.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07805"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "1.841.81"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.084"
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.50"
$ws.Range("E14").Value = "  -3.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6810"
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.448"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008344"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "28.964.61"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.76"
$ws.Range("E19").Value = "  -4.06%  "
$ws.Range("D20").Value = "2.085.13"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.74"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.479"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9997"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.81"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1472"
$ws.Range("E26").Value = "  -5.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.809"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.23"
$ws.Range("E28").Value = "  -3.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.556"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.220"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.166"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.172"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("E33").Value = "  -3.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7692"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.848"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.145"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.679"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01849"
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("D39").Value = "1.238.62"
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9299"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.68"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.812"
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.579"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("E46").Value = "  -3.74%  "
$ws.Range("D47").Value = "1.983.41"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5170"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.56"
$ws.Range("E49").Value = "  -9.35%  "
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.937"
$ws.Range("E51").Value = "  -2.00%  "
